# fix error codes in epp-16-ns-data
# Swap the errorCode values in the DataProvider table (column F) for rows
# 10-14 so that rows which previously held EPP_UNEXPECTED_COMMAND_FAILURE
# now hold EPP_UNEXPECTED_COMMAND_SUCCESS and vice versa.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F10").Value = "EPP_UNEXPECTED_COMMAND_SUCCESS"
$ws.Range("F11").Value = "EPP_UNEXPECTED_COMMAND_FAILURE"
$ws.Range("F12").Value = "EPP_UNEXPECTED_COMMAND_SUCCESS"
$ws.Range("F13").Value = "EPP_UNEXPECTED_COMMAND_FAILURE"
$ws.Range("F14").Value = "EPP_UNEXPECTED_COMMAND_SUCCESS"
